$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 508
$ws.Range("F5").Value = 2351
$ws.Range("F7").Value = 69
$ws.Range("F16").Value = 801
$ws.Range("F17").Value = 3
$ws.Range("F18").Value = 183
$ws.Range("F20").Value = 7459
$ws.Range("F21").Value = 8392
$ws.Range("F34").Value = 1489
$ws.Range("F37").Value = 24
$ws.Range("F39").Value = 30
$ws.Range("F40").Value = 778
$ws.Range("F48").Value = 187
$ws.Range("F49").Value = 26

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 7
$ws.Range("F18").Value = 310

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2644
$ws.Range("F4").Value = 297
$ws.Range("F5").Value = 152

# Sheet "全部类型" (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 152
$ws.Range("F7").Value = 508
$ws.Range("F8").Value = 2351
$ws.Range("F10").Value = 69
$ws.Range("F19").Value = 801
$ws.Range("F21").Value = 183
$ws.Range("F23").Value = 7459
$ws.Range("F24").Value = 7459
$ws.Range("F25").Value = 8392
$ws.Range("F35").Value = 24
$ws.Range("F38").Value = 30
$ws.Range("F41").Value = 778
$ws.Range("F42").Value = 7
$ws.Range("F49").Value = 187
$ws.Range("F50").Value = 310
$ws.Range("F51").Value = 26

$wb.Save()
